# Appends two new daily rows (2025-10-19 / serial 45949) for the two
# charging stations, continuing the existing date/station pattern found
# in the sheet (each date has one row for station "四方坪站充电量(kw)"
# followed by one row for station "高岭站充电量(kw)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 98 : 2025-10-19, 四方坪站充电量(kw) ----
$ws.Range("A98").Value = 45949
$ws.Range("B98").Value = "四方坪站充电量(kw)"
$ws.Range("C98").Value = 580.31700000000001
$ws.Range("D98").Value = 1176.8040000000003
$ws.Range("E98").Value = 520.91800000000012
$ws.Range("F98").Value = 263.53699999999998
$ws.Range("G98").Value = 480.93099999999993
$ws.Range("H98").Value = 575.10699999999997
$ws.Range("I98").Value = 433.06100000000004
$ws.Range("J98").Value = 320.76199999999994
$ws.Range("K98").Value = 134.96799999999999
$ws.Range("L98").Value = 97.206999999999994
$ws.Range("M98").Value = 295.06400000000002
$ws.Range("N98").Value = 179.67999999999998
$ws.Range("O98").Value = 586.91500000000008
$ws.Range("P98").Value = 956.95300000000009
$ws.Range("Q98").Value = 629.37100000000009
$ws.Range("R98").Value = 493.46599999999995
$ws.Range("S98").Value = 279.16799999999995
$ws.Range("T98").Value = 344.44
$ws.Range("U98").Value = 173.32700000000003
$ws.Range("V98").Value = 99.399999999999991
$ws.Range("W98").Value = 97.88000000000001
$ws.Range("X98").Value = 84.35
$ws.Range("Y98").Value = 93.799999999999983
$ws.Range("Z98").Value = 23.62

# ---- Row 99 : 2025-10-19, 高岭站充电量(kw) ----
$ws.Range("A99").Value = 45949
$ws.Range("B99").Value = "高岭站充电量(kw)"
$ws.Range("C99").Value = 329.83699999999999
$ws.Range("D99").Value = 690.73299999999995
$ws.Range("E99").Value = 212.85399999999998
$ws.Range("F99").Value = 82.168999999999997
$ws.Range("G99").Value = 159.64599999999999
$ws.Range("H99").Value = 169.34399999999999
$ws.Range("I99").Value = 185.70400000000001
$ws.Range("J99").Value = 174.73100000000002
$ws.Range("K99").Value = 160.88
$ws.Range("L99").Value = 95.643000000000001
$ws.Range("M99").Value = 75.251999999999995
$ws.Range("N99").Value = 139.10500000000002
$ws.Range("O99").Value = 203.863
$ws.Range("P99").Value = 372.50700000000001
$ws.Range("Q99").Value = 119.598
$ws.Range("R99").Value = 449.84299999999996
$ws.Range("S99").Value = 193.64399999999998
$ws.Range("T99").Value = 84.968000000000004
$ws.Range("U99").Value = 37.36
$ws.Range("V99").Value = 60.201000000000001
$ws.Range("W99").Value = 148.42600000000002
$ws.Range("X99").Value = 56.063000000000002
$ws.Range("Y99").Value = 0
$ws.Range("Z99").Value = 67.164000000000001

# Match the author's final selection (S104) recorded in the sheet view.
$ws.Range("S104").Select() | Out-Null
